$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (I1, J1). Copy the format from the existing
# header cell H1 (bold / centered / bordered style) so the new headers
# re-use the same cell style instead of creating a new one, then set text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I0 / IF numeric columns for rows 2-9.
$data = @{
    2 = @(1, 1)
    3 = @(4, 7)
    4 = @(5, 7)
    5 = @(1, 4)
    6 = @(1, 5)
    7 = @(6, 8)
    8 = @(4, 5)
    9 = @(5, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
